# Scheduled Sheets runner update: refresh Universalis market-board price
# snapshots (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H:N) for a
# batch of leve rows across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1606.5
$ws.Range("I98").Value = 1590.5555
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 1590.5555
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = -92.55549999999994
$ws.Range("N98").Value = -4746
$ws.Range("H105").Value = 26000
$ws.Range("J105").Value = 26000
$ws.Range("L105").Value = 26000
$ws.Range("N105").Value = -32988
$ws.Range("H107").Value = 1879.619
$ws.Range("I107").Value = 1922.2106
$ws.Range("J107").Value = 1475
$ws.Range("K107").Value = 1922.2106
$ws.Range("L107").Value = 1475
$ws.Range("M107").Value = -2.210600000000113
$ws.Range("N107").Value = -5315
$ws.Range("H121").Value = 6945
$ws.Range("J121").Value = 6945
$ws.Range("L121").Value = 20835
$ws.Range("N121").Value = -24329
$ws.Range("H122").Value = 1606.5
$ws.Range("I122").Value = 1590.5555
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 4771.666499999999
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -2321.666499999999
$ws.Range("N122").Value = -10150
$ws.Range("H128").Value = 76925
$ws.Range("J128").Value = 76925
$ws.Range("L128").Value = 76925
$ws.Range("N128").Value = -86885
$ws.Range("H132").Value = 2757.1785
$ws.Range("I132").Value = 1328.16
$ws.Range("J132").Value = 14665.667
$ws.Range("K132").Value = 3984.48
$ws.Range("L132").Value = 43997.001
$ws.Range("M132").Value = -1454.48
$ws.Range("N132").Value = -49057.001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2348.625
$ws.Range("I45").Value = 1277.8
$ws.Range("K45").Value = 1277.8
$ws.Range("M45").Value = -900.8
$ws.Range("H74").Value = 4003.8
$ws.Range("I74").Value = 3673.6667
$ws.Range("K74").Value = 3673.6667
$ws.Range("M74").Value = -2799.6667
$ws.Range("H77").Value = 4003.8
$ws.Range("I77").Value = 3673.6667
$ws.Range("K77").Value = 18368.3335
$ws.Range("M77").Value = -14000.3335
$ws.Range("H97").Value = 410.78946
$ws.Range("J97").Value = 550
$ws.Range("L97").Value = 550
$ws.Range("N97").Value = -1542
$ws.Range("H110").Value = 2579.1904
$ws.Range("I110").Value = 1476.3
$ws.Range("K110").Value = 1476.3
$ws.Range("M110").Value = 568.7

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H20").Value = 2086
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H99").Value = 1783
$ws.Range("I99").Value = 1895
$ws.Range("K99").Value = 1895
$ws.Range("M99").Value = -397
$ws.Range("H105").Value = 4110.4165
$ws.Range("I105").Value = 3732.5
$ws.Range("J105").Value = 4639.5
$ws.Range("K105").Value = 3732.5
$ws.Range("L105").Value = 4639.5
$ws.Range("M105").Value = -1985.5
$ws.Range("N105").Value = -8133.5
$ws.Range("H134").Value = 2682.75
$ws.Range("I134").Value = 3531.4
$ws.Range("J134").Value = 1268.3334
$ws.Range("K134").Value = 10594.2
$ws.Range("L134").Value = 3805.0002
$ws.Range("M134").Value = -8059.200000000001
$ws.Range("N134").Value = -8875.0002

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 122340.89
$ws.Range("I22").Value = 125133.375
$ws.Range("K22").Value = 125133.375
$ws.Range("M22").Value = -124783.375
$ws.Range("H39").Value = 3167
$ws.Range("J39").Value = 4200
$ws.Range("L39").Value = 4200
$ws.Range("N39").Value = -4982
$ws.Range("H49").Value = 3167
$ws.Range("J49").Value = 4200
$ws.Range("L49").Value = 4200
$ws.Range("N49").Value = -4564
$ws.Range("H86").Value = 3300
$ws.Range("I86").Value = 3200
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 3200
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -2077
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 3300
$ws.Range("I89").Value = 3200
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 16000
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -10384
$ws.Range("N89").Value = -28732
$ws.Range("H99").Value = 10523.6875
$ws.Range("I99").Value = 7280.1333
$ws.Range("J99").Value = 13385.647
$ws.Range("K99").Value = 7280.1333
$ws.Range("L99").Value = 13385.647
$ws.Range("M99").Value = -5782.1333
$ws.Range("N99").Value = -16381.647
$ws.Range("H105").Value = 1356.4286
$ws.Range("I105").Value = 1415.8334
$ws.Range("K105").Value = 1415.8334
$ws.Range("M105").Value = 331.1666
$ws.Range("H126").Value = 10523.6875
$ws.Range("I126").Value = 7280.1333
$ws.Range("J126").Value = 13385.647
$ws.Range("K126").Value = 21840.3999
$ws.Range("L126").Value = 40156.94100000001
$ws.Range("M126").Value = -19370.3999
$ws.Range("N126").Value = -45096.94100000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 351.33334
$ws.Range("I2").Value = 23.333334
$ws.Range("J2").Value = 1007.3333
$ws.Range("K2").Value = 140.000004
$ws.Range("L2").Value = 6043.9998
$ws.Range("M2").Value = -27.00000399999999
$ws.Range("N2").Value = -6269.9998
$ws.Range("H8").Value = 248.5
$ws.Range("I8").Value = 248.5
$ws.Range("K8").Value = 745.5
$ws.Range("M8").Value = -606.5
$ws.Range("H23").Value = 125134.625
$ws.Range("J23").Value = 125134.625
$ws.Range("L23").Value = 375403.875
$ws.Range("N23").Value = -375873.875
$ws.Range("H38").Value = 125.22222
$ws.Range("I38").Value = 84.333336
$ws.Range("J38").Value = 207
$ws.Range("K38").Value = 253.000008
$ws.Range("L38").Value = 621
$ws.Range("M38").Value = 93.99999199999999
$ws.Range("N38").Value = -1315
$ws.Range("H107").Value = 125722.125
$ws.Range("J107").Value = 125722.125
$ws.Range("L107").Value = 377166.375
$ws.Range("N107").Value = -381006.375
$ws.Range("H124").Value = 1072.25
$ws.Range("I124").Value = 1072.25
$ws.Range("K124").Value = 3216.75
$ws.Range("M124").Value = 1693.25

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5147.8
$ws.Range("I40").Value = 4997.25
$ws.Range("K40").Value = 4997.25
$ws.Range("M40").Value = -4861.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H136").Value = 2372.4285
$ws.Range("I136").Value = 1305
$ws.Range("K136").Value = 3915
$ws.Range("M136").Value = -1365
